$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously-missing values (blank -> number) ---
$ws.Range("C3").Value = 11.2
$ws.Range("E6").Value = -5.7
$ws.Range("D9").Value = -14.5
$ws.Range("D10").Value = -14.7
$ws.Range("E14").Value = -5.4
$ws.Range("E19").Value = -6.5
$ws.Range("C21").Value = 12.7
$ws.Range("E29").Value = -10
$ws.Range("E30").Value = -5.9
$ws.Range("C34").Value = 10.5

# --- Clear values that became missing (number -> blank) ---
$ws.Range("D4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("E25").ClearContents()

# --- Remove rows no longer present in the data (delete bottom-up so row
#     numbers of not-yet-deleted rows stay stable) ---
$ws.Rows("28").Delete()
$ws.Rows("26").Delete()
